$d = $word.ActiveDocument

# The first paragraph originally reads:
#   "This is a Microsoft word document."
# The edit appends " (Changed main)" right after the existing sentence,
# turning it into:
#   "This is a Microsoft word document. (Changed main)"

$p = $d.Paragraphs.First
$r = $p.Range
$r.InsertAfter(" (")
$r.InsertAfter("Changed main")
$r.InsertAfter(")")
